$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet carries a stray fully-formatted "ghost" row at the very bottom of
# the grid (row 1048576). Drop it so the used range / dimension settles back
# down cleanly once the new content below has been added.
$ws.Rows.Item(1048576).Delete()

# Insert a new row at row 3 for "Transfer line complex angle", pushing the
# Momentum-acceptance row (and everything below it) down by one row.
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = "Transfer line complex angle"
$ws.Cells.Item(3, 2).Value = "TrfLineCmplxAng"
$ws.Cells.Item(3, 3).Value = 8.5
$ws.Cells.Item(3, 4).Value = "°"
$ws.Cells.Item(3, 5).Value = "nuSIM-2021-01"

# Match the row height used by its sibling parameter rows (the header rows
# 1-2 are taller; every data row from here down uses the shorter height).
$ws.Rows.Item(3).RowHeight = 15

# Touch the format of the sheet's new bottom-right corner so the worksheet's
# used range grows to cover through row 11 (rows 10 & 11 stay blank, as in
# the source row layout before this edit).
$ws.Cells.Item(11, 5).NumberFormat = "General"

[void]$ws.Range("F3").Select()
